# repull data, push all data, mean calculation
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update dSF column (F) values that shifted after repulling data
$ws.Range("F10").Value = 0
$ws.Range("F17").Value = -5
$ws.Range("F19").Value = -6
$ws.Range("F21").Value = -1
$ws.Range("F22").Value = -2
$ws.Range("F23").Value = 1
$ws.Range("F28").Value = 1
$ws.Range("F30").Value = 1
$ws.Range("F34").Value = -1
$ws.Range("F35").Value = -1
$ws.Range("F39").Value = 1
$ws.Range("F43").Value = -1

# Update dS0 column (E) value for row 45
$ws.Range("E45").Value = 2
